$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) cost_sheet : rows 2-28, columns H, I, L, M
#    H/L (remanufacturing_*): 600000 -> 9999999
#    I/M (recycling_*):        1000 -> 11000, plus center-aligned style (s="7")
# ---------------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("cost_sheet")

$wsCost.Range("H2:H28").Value = 9999999
$wsCost.Range("L2:L28").Value = 9999999

$wsCost.Range("I2:I28").Value = 11000
$wsCost.Range("I2:I28").HorizontalAlignment = -4108   # xlCenter -> matches style index 7
$wsCost.Range("M2:M28").Value = 11000
$wsCost.Range("M2:M28").HorizontalAlignment = -4108   # xlCenter -> matches style index 7

# widen a few columns (closest achievable widths given the host's column-width
# quantization step)
$wsCost.Columns.Item(6).ColumnWidth = 39
$wsCost.Columns.Item(7).ColumnWidth = 20.666666666666668
$wsCost.Columns.Item(10).ColumnWidth = 23.333333333333332
$wsCost.Columns.Item(11).ColumnWidth = 26.833333333333332

# ---------------------------------------------------------------------------
# 2) Technologies : J3 and J4 (IR EU Secondary) 0.35 -> 0.3888
# ---------------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("Technologies")
$wsTech.Range("J3").Value = 0.38879999999999998
$wsTech.Range("J4").Value = 0.38879999999999998

# ---------------------------------------------------------------------------
# 3) installable_capacity : columns C, D
#    Rows 2-7:   C 13333 -> 25000 , D 23666 -> 50000
#    Rows 8-17:  C  8900 -> 30000 , D 28700 -> 50000
#    Rows 18-25: C  2900 -> 40000 , D 19100 -> 35000
#    Rows 26-28: C  2900 -> 40000 , D 19100 -> 9999999
#    D loses its wrap-text style (becomes plain / unstyled)
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("installable_capacity")

$wsCap.Range("C2:C7").Value = 25000
$wsCap.Range("D2:D7").Value = 50000

$wsCap.Range("C8:C17").Value = 30000
$wsCap.Range("D8:D17").Value = 50000

$wsCap.Range("C18:C25").Value = 40000
$wsCap.Range("D18:D25").Value = 35000

$wsCap.Range("C26:C28").Value = 40000
$wsCap.Range("D26:D28").Value = 9999999

$wsCap.Range("D2:D28").ClearFormats()

# ---------------------------------------------------------------------------
# 4) Selections / active sheet ordering
#    Final state: cost_sheet selection -> F35, Technologies selection -> C4,
#    installable_capacity selection -> C18:C28 and is the active/selected tab.
# ---------------------------------------------------------------------------
$wsCost.Activate()
$wsCost.Range("F35").Select()

$wsTech.Activate()
$wsTech.Range("C4").Select()

$wsCap.Activate()
$wsCap.Range("C18:C28").Select()
